$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BR: set width to match the other data columns (stored width 12) ---
$ws.Range("BR1:BR53").ColumnWidth = 67/6

# --- Header cell BR1: text date label "2024/11/17" (not an Excel date value) ---
$ws.Range("BR1").NumberFormat = "@"
$ws.Range("BR1").Value = "2024/11/17"
$ws.Range("BQ1").Copy()
$ws.Range("BR1").PasteSpecial(-4122)

# --- Data rows 2-53: composite-probability values with conditional highlight styles ---
$ws.Range("N2").Copy()
$ws.Range("BR2").PasteSpecial(-4122)
$ws.Range("BR2").Value = 129.6
$ws.Range("N2").Copy()
$ws.Range("BR3").PasteSpecial(-4122)
$ws.Range("BR3").Value = 139.1
$ws.Range("N2").Copy()
$ws.Range("BR4").PasteSpecial(-4122)
$ws.Range("BR4").Value = 133.3
$ws.Range("A2").Copy()
$ws.Range("BR5").PasteSpecial(-4122)
$ws.Range("BR5").Value = 169.4
$ws.Range("A2").Copy()
$ws.Range("BR6").PasteSpecial(-4122)
$ws.Range("BR6").Value = 145.3
$ws.Range("A2").Copy()
$ws.Range("BR7").PasteSpecial(-4122)
$ws.Range("BR7").Value = 144.9
$ws.Range("A2").Copy()
$ws.Range("BR8").PasteSpecial(-4122)
$ws.Range("BR8").Value = 149
$ws.Range("D2").Copy()
$ws.Range("BR9").PasteSpecial(-4122)
$ws.Range("BR9").Value = 117.6
$ws.Range("A2").Copy()
$ws.Range("BR10").PasteSpecial(-4122)
$ws.Range("BR10").Value = 163
$ws.Range("N2").Copy()
$ws.Range("BR11").PasteSpecial(-4122)
$ws.Range("BR11").Value = 129.4
$ws.Range("D2").Copy()
$ws.Range("BR12").PasteSpecial(-4122)
$ws.Range("BR12").Value = 105.2
$ws.Range("N2").Copy()
$ws.Range("BR13").PasteSpecial(-4122)
$ws.Range("BR13").Value = 133.5
$ws.Range("N2").Copy()
$ws.Range("BR14").PasteSpecial(-4122)
$ws.Range("BR14").Value = 131.1
$ws.Range("N2").Copy()
$ws.Range("BR15").PasteSpecial(-4122)
$ws.Range("BR15").Value = 137.6
$ws.Range("D2").Copy()
$ws.Range("BR16").PasteSpecial(-4122)
$ws.Range("BR16").Value = 116.6
$ws.Range("A2").Copy()
$ws.Range("BR17").PasteSpecial(-4122)
$ws.Range("BR17").Value = 159.3
$ws.Range("A2").Copy()
$ws.Range("BR18").PasteSpecial(-4122)
$ws.Range("BR18").Value = 165.7
$ws.Range("N2").Copy()
$ws.Range("BR19").PasteSpecial(-4122)
$ws.Range("BR19").Value = 129.9
$ws.Range("A2").Copy()
$ws.Range("BR20").PasteSpecial(-4122)
$ws.Range("BR20").Value = 148.3
$ws.Range("A2").Copy()
$ws.Range("BR21").PasteSpecial(-4122)
$ws.Range("BR21").Value = 143.6
$ws.Range("A2").Copy()
$ws.Range("BR22").PasteSpecial(-4122)
$ws.Range("BR22").Value = 168.7
$ws.Range("A2").Copy()
$ws.Range("BR23").PasteSpecial(-4122)
$ws.Range("BR23").Value = 182.5
$ws.Range("N2").Copy()
$ws.Range("BR24").PasteSpecial(-4122)
$ws.Range("BR24").Value = 136
$ws.Range("D2").Copy()
$ws.Range("BR25").PasteSpecial(-4122)
$ws.Range("BR25").Value = 121.8
$ws.Range("N2").Copy()
$ws.Range("BR26").PasteSpecial(-4122)
$ws.Range("BR26").Value = 134.4
$ws.Range("N2").Copy()
$ws.Range("BR27").PasteSpecial(-4122)
$ws.Range("BR27").Value = 131.8
$ws.Range("D2").Copy()
$ws.Range("BR28").PasteSpecial(-4122)
$ws.Range("BR28").Value = 121.5
$ws.Range("A2").Copy()
$ws.Range("BR29").PasteSpecial(-4122)
$ws.Range("BR29").Value = 143.5
$ws.Range("N2").Copy()
$ws.Range("BR30").PasteSpecial(-4122)
$ws.Range("BR30").Value = 139.1
$ws.Range("D2").Copy()
$ws.Range("BR31").PasteSpecial(-4122)
$ws.Range("BR31").Value = 108.1
$ws.Range("D2").Copy()
$ws.Range("BR32").PasteSpecial(-4122)
$ws.Range("BR32").Value = 103.3
$ws.Range("A2").Copy()
$ws.Range("BR33").PasteSpecial(-4122)
$ws.Range("BR33").Value = 140.9
$ws.Range("A2").Copy()
$ws.Range("BR34").PasteSpecial(-4122)
$ws.Range("BR34").Value = 167.5
$ws.Range("N2").Copy()
$ws.Range("BR35").PasteSpecial(-4122)
$ws.Range("BR35").Value = 138.9
$ws.Range("A2").Copy()
$ws.Range("BR36").PasteSpecial(-4122)
$ws.Range("BR36").Value = 167.7
$ws.Range("N2").Copy()
$ws.Range("BR37").PasteSpecial(-4122)
$ws.Range("BR37").Value = 125.5
$ws.Range("A2").Copy()
$ws.Range("BR38").PasteSpecial(-4122)
$ws.Range("BR38").Value = 140
$ws.Range("N2").Copy()
$ws.Range("BR39").PasteSpecial(-4122)
$ws.Range("BR39").Value = 131.9
$ws.Range("A2").Copy()
$ws.Range("BR40").PasteSpecial(-4122)
$ws.Range("BR40").Value = 180.5
$ws.Range("A2").Copy()
$ws.Range("BR41").PasteSpecial(-4122)
$ws.Range("BR41").Value = 203.2
$ws.Range("A2").Copy()
$ws.Range("BR42").PasteSpecial(-4122)
$ws.Range("BR42").Value = 181
$ws.Range("A2").Copy()
$ws.Range("BR43").PasteSpecial(-4122)
$ws.Range("BR43").Value = 148.3
$ws.Range("N2").Copy()
$ws.Range("BR44").PasteSpecial(-4122)
$ws.Range("BR44").Value = 132.5
$ws.Range("A2").Copy()
$ws.Range("BR45").PasteSpecial(-4122)
$ws.Range("BR45").Value = 143
$ws.Range("A2").Copy()
$ws.Range("BR46").PasteSpecial(-4122)
$ws.Range("BR46").Value = 170.7
$ws.Range("A2").Copy()
$ws.Range("BR47").PasteSpecial(-4122)
$ws.Range("BR47").Value = 190.6
$ws.Range("N2").Copy()
$ws.Range("BR48").PasteSpecial(-4122)
$ws.Range("BR48").Value = 138.1
$ws.Range("A2").Copy()
$ws.Range("BR49").PasteSpecial(-4122)
$ws.Range("BR49").Value = 153.9
$ws.Range("N2").Copy()
$ws.Range("BR50").PasteSpecial(-4122)
$ws.Range("BR50").Value = 137
$ws.Range("A2").Copy()
$ws.Range("BR51").PasteSpecial(-4122)
$ws.Range("BR51").Value = 150.9
$ws.Range("A2").Copy()
$ws.Range("BR52").PasteSpecial(-4122)
$ws.Range("BR52").Value = 183.7
$ws.Range("D2").Copy()
$ws.Range("BR53").PasteSpecial(-4122)
$ws.Range("BR53").Value = 113.3

$excel.CutCopyMode = $false
Write-Output "BR column for 2024/11/17 added"
